$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 'nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G4").Value = 'nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G10").Value = 'Omnia.Mohammed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G15").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G16").Value = 'Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range("G17").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G23").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G27").Value = 'marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System, nourhanmohamed@med.asu.edu.eg'
$ws.Range("G29").Value = 'Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G30").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G34").Value = 'nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G35").Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G41").Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G48").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G55").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G58").Value = 'marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System, nourhanmohamed@med.asu.edu.eg'
$ws.Range("G60").Value = 'Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G61").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G64").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G65").Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G66").Value = 'nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G72").Value = 'Omnia.Mohammed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G76").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G79").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G89").Value = 'Aya.alaa-eldein@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G91").Value = 'Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G92").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G95").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G96").Value = 'nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G97").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G103").Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G108").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G109").Value = 'Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range("G116").Value = 'heba.e.mohamed@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G117").Value = 'heba.e.mohamed@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G120").Value = 'dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G122").Value = 'Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G123").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G126").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G127").Value = 'nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G128").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G138").Value = 'Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G139").Value = 'Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G147").Value = 'heba.e.mohamed@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G148").Value = 'heba.e.mohamed@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G151").Value = 'Amr-Saeed@med.asu.edu.eg, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com'
$ws.Range("G153").Value = 'Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G154").Value = 'Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G157").Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G158").Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G159").Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G170").Value = 'Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G179").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G182").Value = 'dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G184").Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nahla91@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G185").Value = 'Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
